$wb = $excel.ActiveWorkbook

# --- DBS sheet (sheet2): append new row 13 with the new CdCode entry ---
$wsDBS = $wb.Worksheets.Item("DBS")

$wsDBS.Range("A13").Value = "defItemEq2"
$wsDBS.Range("B13").Value = "DefCode = ,AND DefType = ,AND Item %"
$wsDBS.Range("C13").Value = "Code asc"

# --- Switch the active/selected tab from DBD to DBS, matching the commit ---
$wsDBS.Activate() | Out-Null
$wsDBS.Range("B17").Select() | Out-Null

Write-Output "done"
